$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells store numeric-looking text (e.g. "58.442.86",
# "1.00", using "." as both decimal and thousands separator per the source
# feed). A plain Range.Value assignment of a numeric-looking string gets
# auto-coerced by Excel into a real Number (dropping trailing zeros / mangling
# the thousands-grouped values), so each such cell is forced to Text via
# NumberFormat "@" before the value is written, then the format is cleared
# again afterwards so the cell keeps its original (default, unstyled) look -
# only the stored value differs from the source workbook.
# (NumberFormat / ClearFormats must be applied per-cell: applying them to a
# multi-area Union range only affects the first area.)
$priceRefs = @("D2", "D3", "D5", "D6", "D9", "D13", "D14", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D27", "D32", "D33", "D35", "D37", "D38", "D40", "D42", "D43", "D44", "D48", "D50", "D51")
foreach ($ref in $priceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.442.86"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "2.619.27"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "532.99"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "142.21"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "6.91"
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "3.085.80"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "58.374.38"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").Value = "20.63"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "2.625.17"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "334.31"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "10.10"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "66.33"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "7.05"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "18.75"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "150.34"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("D35").Value = "0.854"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "0.808"
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "1.41"
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "279.54"
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "0.592"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").Value = "10.69"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "18.92"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "1.935.98"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "17.84"
$ws.Range("E50").Value = "  -4.21%  "
$ws.Range("D51").Value = "112.94"
$ws.Range("E51").Value = "  +0.82%  "

foreach ($ref in $priceRefs) {
    $ws.Range($ref).ClearFormats()
}
